$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.018172740936279
$ws.Range("B1").Value = 2.293184518814087
$ws.Range("C1").Value = 4.776578426361084
$ws.Range("D1").Value = 1.504894971847534
$ws.Range("E1").Value = 1.274855852127075
